$wb = $excel.ActiveWorkbook

# --- Sheet "RECURSOS NUEVOS" (3rd tab): fill in the "FICHA" (B) column ---
$wsNuevos = $wb.Worksheets.Item("RECURSOS NUEVOS")

$fichaValues = @{
    2  = "no"
    3  = "no"
    4  = "no"
    5  = "no"
    6  = "no"
    7  = "no"
    8  = "sí"
    9  = "sí"
    10 = "sí"
    11 = "no"
    12 = "no"
    13 = "no"
    14 = "sí"
    15 = "no"
    16 = "no"
    17 = "no"
    18 = "no"
}

$rowsNeedingFormatClear = @(17, 18)

foreach ($row in $fichaValues.Keys) {
    $cell = $wsNuevos.Cells.Item($row, 2)
    if ($rowsNeedingFormatClear -contains $row) {
        $cell.ClearFormats()
    }
    $cell.Value = $fichaValues[$row]
}

# Update the active selection on that sheet (was A23, now A20)
$wsNuevos.Range("A20").Select() | Out-Null

# --- Sheet "GUION" (1st tab) becomes the selected / active tab ---
$wsGuion = $wb.Worksheets.Item("GUION")
$wsGuion.Activate()
$wsGuion.Range("B7").Select() | Out-Null
